$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Full updated data set for rows 2-10 (TPM re-run added an "Inflammatory-Mac"
# target cluster and recomputed all the expression / specificity figures).
# Columns A:T match the header row:
# A Sending cluster, B Ligand symbol, C Receptor symbol, D Target cluster,
# E Ligand-expressing cells ... T Edge total expression derived specificity
$data = @(
    @("ECs","Wnt9a","Fzd10","FAPs",3,1,0.6558959999999999,1.967688,0.1303533075142512,0.1303533075142512,1,0.3333333333333333,0.005445666666666667,0.016337,0.1242225162340131,0.1242225162340131,0.003571790984,0.032146118856,0.01619281585884638,0.01619281585884638),
    @("ECs","Wnt9a","Fzd10","Inflammatory-Mac",3,1,0.6558959999999999,1.967688,0.1303533075142512,0.1303533075142512,2,0.6666666666666666,0.03095266666666667,0.092858,0.706069315814286,0.7060693158142859,0.020301730256,0.182715572304,0.09203847065071658,0.09203847065071656),
    @("ECs","Wnt9a","Fzd10","MuSCs",3,1,0.6558959999999999,1.967688,0.1303533075142512,0.1303533075142512,1,0.3333333333333333,0.007439666666666666,0.022319,0.169708167951701,0.169708167951701,0.004879647607999999,0.04391682847199999,0.02212202100468827,0.02212202100468827),
    @("FAPs","Wnt9a","Fzd10","FAPs",3,1,3.335609,10.006827,0.6629216609406124,0.6629216609406123,1,0.3333333333333333,0.005445666666666667,0.016337,0.1242225162340131,0.1242225162340131,0.01816461474433333,0.163481532699,0.08234979678807418,0.08234979678807415),
    @("FAPs","Wnt9a","Fzd10","Inflammatory-Mac",3,1,3.335609,10.006827,0.6629216609406124,0.6629216609406123,2,0.6666666666666666,0.03095266666666667,0.092858,0.706069315814286,0.7060693158142859,0.1032459935073333,0.9292139415659999,0.4680686435788083,0.4680686435788082),
    @("FAPs","Wnt9a","Fzd10","MuSCs",3,1,3.335609,10.006827,0.6629216609406124,0.6629216609406123,1,0.3333333333333333,0.007439666666666666,0.022319,0.169708167951701,0.169708167951701,0.02481581909033333,0.223342371813,0.11250322057373,0.11250322057373),
    @("MuSCs","Wnt9a","Fzd10","FAPs",3,1,1.040174,3.120522,0.2067250315451364,0.2067250315451363,1,0.3333333333333333,0.005445666666666667,0.016337,0.1242225162340131,0.1242225162340131,0.005664440879333335,0.05097996791400001,0.02567990358709258,0.02567990358709257),
    @("MuSCs","Wnt9a","Fzd10","Inflammatory-Mac",3,1,1.040174,3.120522,0.2067250315451364,0.2067250315451363,2,0.6666666666666666,0.03095266666666667,0.092858,0.706069315814286,0.7060693158142859,0.03219615909733334,0.289765431876,0.1459622015847611,0.1459622015847611),
    @("MuSCs","Wnt9a","Fzd10","MuSCs",3,1,1.040174,3.120522,0.2067250315451364,0.2067250315451363,1,0.3333333333333333,0.007439666666666666,0.022319,0.169708167951701,0.169708167951701,0.007738547835333334,0.069646930518,0.03508292637328268,0.03508292637328268)
)

$startRow = 2
for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $data[$i]
    for ($j = 0; $j -lt $row.Length; $j++) {
        $ws.Cells.Item($startRow + $i, $j + 1).Value = $row[$j]
    }
}
